$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the first (and only) medication record, replacing the old data.
# A2 and D2 must stay text (not auto-converted to numbers), matching the
# original workbook's inlineStr typing for those columns.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value = "Aciclovir 200mg"
$ws.Range("C2").Value = "CIMED INDUSTRIA S.A"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "143810181"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "200 MG COM CT BL AL PLAS TRANS X 10  Ativo"
# F2 ("OK") is unchanged.

# Remove the now-obsolete rows 3-5 entirely, shrinking the used range to A1:F2.
$ws.Range("A3:F5").Delete()
